$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 29; this shifts the existing rows 29-39
# down to rows 32-42 (all their content/styles move with them).
$ws.Range("A29:A31").EntireRow.Insert()

# Common values shared by every row in this data block.
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria = "Damasco"

# New row 29: Modesto / Especial
$r = 29
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44579
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Modesto"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 21000
$ws.Cells.Item($r, 15).Value = 21000
$ws.Cells.Item($r, 16).Value = 21000
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región Metropolitana"
$ws.Cells.Item($r, 19).Value = 1167
$ws.Cells.Item($r, 20).Value = 18

# New row 30: Modesto / Primera
$r = 30
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44579
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Modesto"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 18000
$ws.Cells.Item($r, 16).Value = 18000
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región Metropolitana"
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 18

# New row 31: Modesto / Segunda
$r = 31
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44579
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Modesto"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 16000
$ws.Cells.Item($r, 15).Value = 16000
$ws.Cells.Item($r, 16).Value = 16000
$ws.Cells.Item($r, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item($r, 18).Value = "Región Metropolitana"
$ws.Cells.Item($r, 19).Value = 889
$ws.Cells.Item($r, 20).Value = 18
